$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4200
$ws.Range("L3").Value = 4451
$ws.Range("L4").Value = 1095
$ws.Range("L6").Value = 3847
$ws.Range("L7").Value = 13847
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L3").Value = 42
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 161
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 264
$ws.Range("L3").Value = 306
$ws.Range("L6").Value = 250
$ws.Range("L7").Value = 917
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 178
$ws.Range("L3").Value = 217
$ws.Range("L4").Value = 37
$ws.Range("L7").Value = 646
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 151
$ws.Range("L6").Value = 139
$ws.Range("L7").Value = 511
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 92
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 264
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L4").Value = 20
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 232
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 460
$ws.Range("L8").Value = 917
$ws.Range("L10").Value = 91
$ws.Range("L11").Value = 221
$ws.Range("L13").Value = 23
$ws.Range("L19").Value = 390
$ws.Range("L20").Value = 348
$ws.Range("L22").Value = 41
$ws.Range("L24").Value = 37
$ws.Range("L29").Value = 774
$ws.Range("L33").Value = 646
$ws.Range("L36").Value = 174
$ws.Range("L37").Value = 511
$ws.Range("L48").Value = 181
$ws.Range("L51").Value = 170
$ws.Range("L52").Value = 280
$ws.Range("L53").Value = 161
$ws.Range("L55").Value = 134
$ws.Range("L63").Value = 43
$ws.Range("L65").Value = 264
$ws.Range("L67").Value = 474
$ws.Range("L68").Value = 42
$ws.Range("L78").Value = 178
$ws.Range("L79").Value = 363
$ws.Range("L84").Value = 136
$ws.Range("L85").Value = 721
$ws.Range("L88").Value = 151
$ws.Range("L89").Value = 197
$ws.Range("L90").Value = 140
$ws.Range("L91").Value = 194
$ws.Range("L93").Value = 72
$ws.Range("L94").Value = 174
$ws.Range("L96").Value = 147
$ws.Range("L97").Value = 116
$ws.Range("L98").Value = 79
$ws.Range("L99").Value = 232
$ws.Range("L101").Value = 13847
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 134
$ws.Range("L6").Value = 113
$ws.Range("L7").Value = 474
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 136
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 234
$ws.Range("L3").Value = 293
$ws.Range("L6").Value = 199
$ws.Range("L7").Value = 774
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 79
$ws.Range("L7").Value = 181
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 119
$ws.Range("L7").Value = 390
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 139
$ws.Range("L6").Value = 123
$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L5").Value = 13
$ws.Range("L6").Value = 23
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L2").Value = 39
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 91
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 55
$ws.Range("L7").Value = 178
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 43
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 134
$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L3").Value = 12
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 37
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L4").Value = 18
$ws.Range("L7").Value = 147
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 81
$ws.Range("L7").Value = 194
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 122
$ws.Range("L3").Value = 130
$ws.Range("L7").Value = 363
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 108
$ws.Range("L7").Value = 348
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 174
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 72
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 153
$ws.Range("L6").Value = 120
$ws.Range("L7").Value = 460
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 69
$ws.Range("L7").Value = 174
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L3").Value = 14
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 79
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 221
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L6").Value = 60
$ws.Range("L7").Value = 116
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 42
$ws.Range("L7").Value = 151
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L6").Value = 56
$ws.Range("L7").Value = 197
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 140
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L4").Value = 25
$ws.Range("L7").Value = 170
$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 42
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 220
$ws.Range("L3").Value = 290
$ws.Range("L4").Value = 47
$ws.Range("L7").Value = 721
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 41
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 280
